# Auto-generated edit script: refresh market-price derived values
# across 8 item-leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$edits = @(
    @{Sheet="ALC"; Cell="H4"; Value=273.17648},
    @{Sheet="ALC"; Cell="I4"; Value=126.92308},
    @{Sheet="ALC"; Cell="K4"; Value=126.92308},
    @{Sheet="ALC"; Cell="M4"; Value=-12.92308},
    @{Sheet="ALC"; Cell="H43"; Value=1997.5},
    @{Sheet="ALC"; Cell="I43"; Value=1990},
    @{Sheet="ALC"; Cell="K43"; Value=1990},
    @{Sheet="ALC"; Cell="M43"; Value=-1921},
    @{Sheet="ALC"; Cell="H70"; Value=1112.5},
    @{Sheet="ALC"; Cell="I70"; Value=1135},
    @{Sheet="ALC"; Cell="J70"; Value=1000},
    @{Sheet="ALC"; Cell="K70"; Value=3405},
    @{Sheet="ALC"; Cell="L70"; Value=3000},
    @{Sheet="ALC"; Cell="M70"; Value=-3135},
    @{Sheet="ALC"; Cell="N70"; Value=-3540},
    @{Sheet="ALC"; Cell="H73"; Value=1112.5},
    @{Sheet="ALC"; Cell="I73"; Value=1135},
    @{Sheet="ALC"; Cell="J73"; Value=1000},
    @{Sheet="ALC"; Cell="K73"; Value=3405},
    @{Sheet="ALC"; Cell="L73"; Value=3000},
    @{Sheet="ALC"; Cell="M73"; Value=-2469},
    @{Sheet="ALC"; Cell="N73"; Value=-4872},
    @{Sheet="ALC"; Cell="H92"; Value=919.2857},
    @{Sheet="ALC"; Cell="I92"; Value=954},
    @{Sheet="ALC"; Cell="J92"; Value=832.5},
    @{Sheet="ALC"; Cell="K92"; Value=954},
    @{Sheet="ALC"; Cell="L92"; Value=832.5},
    @{Sheet="ALC"; Cell="M92"; Value=294},
    @{Sheet="ALC"; Cell="N92"; Value=-3328.5},
    @{Sheet="ALC"; Cell="H129"; Value=132367.78},
    @{Sheet="ALC"; Cell="J129"; Value=147884.23},
    @{Sheet="ALC"; Cell="L129"; Value=443652.6900000001},
    @{Sheet="ALC"; Cell="N129"; Value=-453652.6900000001},
    @{Sheet="ALC"; Cell="H132"; Value=2942.8},
    @{Sheet="ALC"; Cell="I132"; Value=3122.0417},
    @{Sheet="ALC"; Cell="J132"; Value=2225.8333},
    @{Sheet="ALC"; Cell="K132"; Value=9366.125100000001},
    @{Sheet="ALC"; Cell="L132"; Value=6677.499899999999},
    @{Sheet="ALC"; Cell="M132"; Value=-6836.125100000001},
    @{Sheet="ALC"; Cell="N132"; Value=-11737.4999},
    @{Sheet="ALC"; Cell="H138"; Value=1598.4103},
    @{Sheet="ALC"; Cell="J138"; Value=2681.8948},
    @{Sheet="ALC"; Cell="L138"; Value=8045.6844},
    @{Sheet="ALC"; Cell="N138"; Value=-18325.6844},
    @{Sheet="ARM"; Cell="H21"; Value=0},
    @{Sheet="ARM"; Cell="I21"; Value=0},
    @{Sheet="ARM"; Cell="K21"; Value=0},
    @{Sheet="ARM"; Cell="M21"; Value=$null},
    @{Sheet="ARM"; Cell="H32"; Value=2553.18},
    @{Sheet="ARM"; Cell="I32"; Value=2327.6904},
    @{Sheet="ARM"; Cell="J32"; Value=3737},
    @{Sheet="ARM"; Cell="K32"; Value=2327.6904},
    @{Sheet="ARM"; Cell="L32"; Value=3737},
    @{Sheet="ARM"; Cell="M32"; Value=-2040.6904},
    @{Sheet="ARM"; Cell="N32"; Value=-4311},
    @{Sheet="ARM"; Cell="H74"; Value=2590.1777},
    @{Sheet="ARM"; Cell="I74"; Value=2736.7104},
    @{Sheet="ARM"; Cell="J74"; Value=1794.7142},
    @{Sheet="ARM"; Cell="K74"; Value=2736.7104},
    @{Sheet="ARM"; Cell="L74"; Value=1794.7142},
    @{Sheet="ARM"; Cell="M74"; Value=-1862.7104},
    @{Sheet="ARM"; Cell="N74"; Value=-3542.7142},
    @{Sheet="ARM"; Cell="H77"; Value=2590.1777},
    @{Sheet="ARM"; Cell="I77"; Value=2736.7104},
    @{Sheet="ARM"; Cell="J77"; Value=1794.7142},
    @{Sheet="ARM"; Cell="K77"; Value=13683.552},
    @{Sheet="ARM"; Cell="L77"; Value=8973.571},
    @{Sheet="ARM"; Cell="M77"; Value=-9315.552},
    @{Sheet="ARM"; Cell="N77"; Value=-17709.571},
    @{Sheet="ARM"; Cell="H102"; Value=4800},
    @{Sheet="ARM"; Cell="I102"; Value=1466.6666},
    @{Sheet="ARM"; Cell="J102"; Value=6800},
    @{Sheet="ARM"; Cell="K102"; Value=1466.6666},
    @{Sheet="ARM"; Cell="L102"; Value=6800},
    @{Sheet="ARM"; Cell="M102"; Value=155.3334},
    @{Sheet="ARM"; Cell="N102"; Value=-10044},
    @{Sheet="ARM"; Cell="H132"; Value=20924.334},
    @{Sheet="ARM"; Cell="I132"; Value=2134.5},
    @{Sheet="ARM"; Cell="J132"; Value=103599.6},
    @{Sheet="ARM"; Cell="K132"; Value=6403.5},
    @{Sheet="ARM"; Cell="L132"; Value=310798.8},
    @{Sheet="ARM"; Cell="M132"; Value=-3873.5},
    @{Sheet="ARM"; Cell="N132"; Value=-315858.8},
    @{Sheet="BSM"; Cell="H22"; Value=281.2},
    @{Sheet="BSM"; Cell="I22"; Value=284.66666},
    @{Sheet="BSM"; Cell="J22"; Value=250},
    @{Sheet="BSM"; Cell="K22"; Value=284.66666},
    @{Sheet="BSM"; Cell="L22"; Value=250},
    @{Sheet="BSM"; Cell="M22"; Value=-111.66666},
    @{Sheet="BSM"; Cell="N22"; Value=-596},
    @{Sheet="BSM"; Cell="H94"; Value=3385.76},
    @{Sheet="BSM"; Cell="I94"; Value=1574.9333},
    @{Sheet="BSM"; Cell="J94"; Value=6102},
    @{Sheet="BSM"; Cell="K94"; Value=1574.9333},
    @{Sheet="BSM"; Cell="L94"; Value=6102},
    @{Sheet="BSM"; Cell="M94"; Value=-1123.9333},
    @{Sheet="BSM"; Cell="N94"; Value=-7004},
    @{Sheet="BSM"; Cell="H134"; Value=4258},
    @{Sheet="BSM"; Cell="I134"; Value=4582.4707},
    @{Sheet="BSM"; Cell="J134"; Value=1500},
    @{Sheet="BSM"; Cell="K134"; Value=13747.4121},
    @{Sheet="BSM"; Cell="L134"; Value=4500},
    @{Sheet="BSM"; Cell="M134"; Value=-11212.4121},
    @{Sheet="BSM"; Cell="N134"; Value=-9570},
    @{Sheet="CRP"; Cell="H31"; Value=3377.3225},
    @{Sheet="CRP"; Cell="I31"; Value=2829.2856},
    @{Sheet="CRP"; Cell="J31"; Value=3828.647},
    @{Sheet="CRP"; Cell="K31"; Value=2829.2856},
    @{Sheet="CRP"; Cell="L31"; Value=3828.647},
    @{Sheet="CRP"; Cell="M31"; Value=-2534.2856},
    @{Sheet="CRP"; Cell="N31"; Value=-4418.647},
    @{Sheet="CRP"; Cell="H34"; Value=3377.3225},
    @{Sheet="CRP"; Cell="I34"; Value=2829.2856},
    @{Sheet="CRP"; Cell="J34"; Value=3828.647},
    @{Sheet="CRP"; Cell="K34"; Value=2829.2856},
    @{Sheet="CRP"; Cell="L34"; Value=3828.647},
    @{Sheet="CRP"; Cell="M34"; Value=-2627.2856},
    @{Sheet="CRP"; Cell="N34"; Value=-4232.647},
    @{Sheet="CRP"; Cell="H99"; Value=18821030},
    @{Sheet="CRP"; Cell="I99"; Value=3971144.5},
    @{Sheet="CRP"; Cell="K99"; Value=3971144.5},
    @{Sheet="CRP"; Cell="M99"; Value=-3969646.5},
    @{Sheet="CRP"; Cell="H105"; Value=17857688},
    @{Sheet="CRP"; Cell="I105"; Value=17857688},
    @{Sheet="CRP"; Cell="K105"; Value=17857688},
    @{Sheet="CRP"; Cell="M105"; Value=-17855941},
    @{Sheet="CRP"; Cell="H126"; Value=18821030},
    @{Sheet="CRP"; Cell="I126"; Value=3971144.5},
    @{Sheet="CRP"; Cell="K126"; Value=11913433.5},
    @{Sheet="CRP"; Cell="M126"; Value=-11910963.5},
    @{Sheet="CUL"; Cell="H36"; Value=2389.5},
    @{Sheet="CUL"; Cell="I36"; Value=1682.75},
    @{Sheet="CUL"; Cell="K36"; Value=5048.25},
    @{Sheet="CUL"; Cell="M36"; Value=-4879.25},
    @{Sheet="CUL"; Cell="H75"; Value=1534.75},
    @{Sheet="CUL"; Cell="I75"; Value=963},
    @{Sheet="CUL"; Cell="J75"; Value=3250},
    @{Sheet="CUL"; Cell="K75"; Value=2889},
    @{Sheet="CUL"; Cell="L75"; Value=9750},
    @{Sheet="CUL"; Cell="M75"; Value=-1891},
    @{Sheet="CUL"; Cell="N75"; Value=-11746},
    @{Sheet="CUL"; Cell="H78"; Value=1534.75},
    @{Sheet="CUL"; Cell="I78"; Value=963},
    @{Sheet="CUL"; Cell="J78"; Value=3250},
    @{Sheet="CUL"; Cell="K78"; Value=8667},
    @{Sheet="CUL"; Cell="L78"; Value=29250},
    @{Sheet="CUL"; Cell="M78"; Value=-3675},
    @{Sheet="CUL"; Cell="N78"; Value=-39234},
    @{Sheet="CUL"; Cell="H131"; Value=802.81055},
    @{Sheet="CUL"; Cell="I131"; Value=564.75},
    @{Sheet="CUL"; Cell="J131"; Value=813.2747000000001},
    @{Sheet="CUL"; Cell="K131"; Value=1694.25},
    @{Sheet="CUL"; Cell="L131"; Value=2439.8241},
    @{Sheet="CUL"; Cell="M131"; Value=3345.75},
    @{Sheet="CUL"; Cell="N131"; Value=-12519.8241},
    @{Sheet="CUL"; Cell="H137"; Value=2224.3333},
    @{Sheet="CUL"; Cell="I137"; Value=671.6667},
    @{Sheet="CUL"; Cell="J137"; Value=2446.1428},
    @{Sheet="CUL"; Cell="K137"; Value=2015.0001},
    @{Sheet="CUL"; Cell="L137"; Value=7338.428400000001},
    @{Sheet="CUL"; Cell="M137"; Value=3084.9999},
    @{Sheet="CUL"; Cell="N137"; Value=-17538.4284},
    @{Sheet="GSM"; Cell="H22"; Value=100},
    @{Sheet="GSM"; Cell="J22"; Value=100},
    @{Sheet="GSM"; Cell="L22"; Value=100},
    @{Sheet="GSM"; Cell="N22"; Value=-1158},
    @{Sheet="GSM"; Cell="H97"; Value=2885},
    @{Sheet="GSM"; Cell="J97"; Value=6896.6665},
    @{Sheet="GSM"; Cell="L97"; Value=6896.6665},
    @{Sheet="GSM"; Cell="N97"; Value=-7888.6665},
    @{Sheet="GSM"; Cell="H113"; Value=3787.5},
    @{Sheet="GSM"; Cell="I113"; Value=2600},
    @{Sheet="GSM"; Cell="J113"; Value=4975},
    @{Sheet="GSM"; Cell="K113"; Value=2600},
    @{Sheet="GSM"; Cell="L113"; Value=4975},
    @{Sheet="GSM"; Cell="M113"; Value=-430},
    @{Sheet="GSM"; Cell="N113"; Value=-9315},
    @{Sheet="GSM"; Cell="H121"; Value=30000},
    @{Sheet="GSM"; Cell="J121"; Value=30000},
    @{Sheet="GSM"; Cell="L121"; Value=30000},
    @{Sheet="GSM"; Cell="N121"; Value=-33494},
    @{Sheet="GSM"; Cell="H123"; Value=10325.8},
    @{Sheet="GSM"; Cell="I123"; Value=0},
    @{Sheet="GSM"; Cell="J123"; Value=10325.8},
    @{Sheet="GSM"; Cell="K123"; Value=0},
    @{Sheet="GSM"; Cell="L123"; Value=10325.8},
    @{Sheet="GSM"; Cell="M123"; Value=$null},
    @{Sheet="GSM"; Cell="N123"; Value=-15225.8},
    @{Sheet="GSM"; Cell="H132"; Value=61088.668},
    @{Sheet="GSM"; Cell="I132"; Value=6300},
    @{Sheet="GSM"; Cell="K132"; Value=18900},
    @{Sheet="GSM"; Cell="M132"; Value=-16370},
    @{Sheet="LTW"; Cell="H16"; Value=284},
    @{Sheet="LTW"; Cell="I16"; Value=261.42856},
    @{Sheet="LTW"; Cell="J16"; Value=600},
    @{Sheet="LTW"; Cell="K16"; Value=261.42856},
    @{Sheet="LTW"; Cell="L16"; Value=600},
    @{Sheet="LTW"; Cell="M16"; Value=-91.42856},
    @{Sheet="LTW"; Cell="N16"; Value=-940},
    @{Sheet="LTW"; Cell="H33"; Value=1507.5},
    @{Sheet="LTW"; Cell="I33"; Value=1507.5},
    @{Sheet="LTW"; Cell="K33"; Value=1507.5},
    @{Sheet="LTW"; Cell="M33"; Value=-1217.5},
    @{Sheet="LTW"; Cell="H55"; Value=271.75},
    @{Sheet="LTW"; Cell="I55"; Value=209.8},
    @{Sheet="LTW"; Cell="J55"; Value=316},
    @{Sheet="LTW"; Cell="K55"; Value=209.8},
    @{Sheet="LTW"; Cell="L55"; Value=316},
    @{Sheet="LTW"; Cell="M55"; Value=-36.80000000000001},
    @{Sheet="LTW"; Cell="N55"; Value=-662},
    @{Sheet="LTW"; Cell="H119"; Value=32000},
    @{Sheet="LTW"; Cell="J119"; Value=32000},
    @{Sheet="LTW"; Cell="L119"; Value=32000},
    @{Sheet="LTW"; Cell="N119"; Value=-41676},
    @{Sheet="WVR"; Cell="H74"; Value=39909},
    @{Sheet="WVR"; Cell="J74"; Value=39909},
    @{Sheet="WVR"; Cell="L74"; Value=39909},
    @{Sheet="WVR"; Cell="N74"; Value=-41781},
    @{Sheet="WVR"; Cell="H77"; Value=39909},
    @{Sheet="WVR"; Cell="J77"; Value=39909},
    @{Sheet="WVR"; Cell="L77"; Value=119727},
    @{Sheet="WVR"; Cell="N77"; Value=-129087},
    @{Sheet="WVR"; Cell="H81"; Value=1335.75},
    @{Sheet="WVR"; Cell="I81"; Value=1335.75},
    @{Sheet="WVR"; Cell="K81"; Value=2671.5},
    @{Sheet="WVR"; Cell="M81"; Value=-1610.5},
    @{Sheet="WVR"; Cell="H84"; Value=1335.75},
    @{Sheet="WVR"; Cell="I84"; Value=1335.75},
    @{Sheet="WVR"; Cell="K84"; Value=13357.5},
    @{Sheet="WVR"; Cell="M84"; Value=-8053.5}
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    if ($null -eq $e.Value) {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value2 = $e.Value
    }
}
